$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Replace-Text "2024-01-18 Thursday" "2024-01-19 Friday"
Replace-Text "36+48=84" "58-4=54"
Replace-Text "97-52=45" "71+1=72"
Replace-Text "25+15=40" "73-17=56"
Replace-Text "84-42=42" "71-36=35"
Replace-Text "94-85=9" "7+15=22"
Replace-Text "65-58=7" "33-30=3"
Replace-Text "49-9=40" "97-84=13"
Replace-Text "41-24=17" "2+79=81"
Replace-Text "41+56=97" "42-7=35"
Replace-Text "11+62=73" "88-15=73"
Replace-Text "40+0=40" "66-52=14"
Replace-Text "26-18=8" "51-2=49"
Replace-Text "70-20=50" "0+44=44"
Replace-Text "90-40=50" "16+44=60"
Replace-Text "6-6=0" "85-12=73"
Replace-Text "42-19=23" "6-4=2"
Replace-Text "14+5=19" "4+14=18"
Replace-Text "10+54=64" "76-50=26"
Replace-Text "3+33=36" "36-13=23"
Replace-Text "66+28=94" "52-11=41"
Replace-Text "78-58=20" "63-60=3"
Replace-Text "51-24=27" "70-10=60"
Replace-Text "35+63=98" "41+6=47"
Replace-Text "41-8=33" "49-37=12"
Replace-Text "20-19=1" "22+38=60"
Replace-Text "73-6=67" "1+3=4"
Replace-Text "77-20=57" "47+3=50"
Replace-Text "95-62=33" "76+11=87"
Replace-Text "61-6=55" "70-28=42"
Replace-Text "72-43=29" "8+49=57"
Replace-Text "53+38=91" "80-40=40"
Replace-Text "77-41=36" "52-24=28"
Replace-Text "6+79=85" "20+51=71"
Replace-Text "40+1=41" "88-53=35"
Replace-Text "76-37=39" "71+9=80"
Replace-Text "41-7=34" "66+5=71"
Replace-Text "22-15=7" "85-77=8"
Replace-Text "10+8=18" "76-21=55"
Replace-Text "50-11=39" "74-17=57"
Replace-Text "0+36=36" "27+12=39"
Replace-Text "21+77=98" "61-55=6"
Replace-Text "74-14=60" "15+1=16"
Replace-Text "42+38=80" "29+45=74"
Replace-Text "91-62=29" "28+46=74"
Replace-Text "25+61=86" "80-80=0"
Replace-Text "99-58=41" "20+76=96"
Replace-Text "8+91=99" "24+49=73"
Replace-Text "78-28=50" "62-9=53"
Replace-Text "60-36=24" "41-36=5"
Replace-Text "22+16=38" "25+24=49"
Replace-Text "25-19=6" "27+21=48"
Replace-Text "50+46=96" "49+19=68"
Replace-Text "69+7=76" "86-10=76"
Replace-Text "27+39=66" "65-55=10"
Replace-Text "68-18=50" "73+21=94"
Replace-Text "34+46=80" "50-36=14"
Replace-Text "88+0=88" "79-69=10"
Replace-Text "4+95=99" "92-48=44"
Replace-Text "58-10=48" "28+33=61"
Replace-Text "84-83=1" "43+40=83"
Replace-Text "59-47=12" "69-13=56"
Replace-Text "84-78=6" "69-15=54"
Replace-Text "31-6=25" "24+8=32"
Replace-Text "41-13=28" "3+86=89"
Replace-Text "53-20=33" "22-12=10"
Replace-Text "29+6=35" "90-46=44"
Replace-Text "96-14=82" "51-17=34"
Replace-Text "48-2=46" "18+31=49"
Replace-Text "3+64=67" "5+89=94"
Replace-Text "45+0=45" "85+6=91"
Replace-Text "29-18=11" "77+20=97"
Replace-Text "4+38=42" "68-12=56"
Replace-Text "86-43=43" "76-53=23"
Replace-Text "58+21=79" "31+68=99"
Replace-Text "86-51=35" "47-33=14"
Replace-Text "30-11=19" "71-5=66"
Replace-Text "31+15=46" "60-60=0"
Replace-Text "11-8=3" "35+20=55"
Replace-Text "7+33=40" "89-57=32"
Replace-Text "97-25=72" "57-52=5"
Replace-Text "38+30=68" "59-38=21"
Replace-Text "84-22=62" "63+31=94"
Replace-Text "37-31=6" "42-16=26"
Replace-Text "67-53=14" "0+74=74"
Replace-Text "70-32=38" "38-6=32"
Replace-Text "50+14=64" "43+28=71"
Replace-Text "49-28=21" "99-37=62"
Replace-Text "7+77=84" "71-67=4"
Replace-Text "98-52=46" "96-22=74"
Replace-Text "95-54=41" "45-41=4"
Replace-Text "14+84=98" "93-39=54"
Replace-Text "67+28=95" "72+22=94"
Replace-Text "12-5=7" "86-57=29"
Replace-Text "66-59=7" "32-12=20"
Replace-Text "96+0=96" "98-71=27"
Replace-Text "22+25=47" "3+96=99"
Replace-Text "92-27=65" "63+11=74"
Replace-Text "48+30=78" "18+4=22"
Replace-Text "54+39=93" "0+82=82"
Replace-Text "18+79=97" "25-13=12"
